# Edges.pptx edit script
# - Bumps the cached "datetimeFigureOut" footer date text on the slide
#   master and every slide layout from 2024-07-12 to 2024-07-20.
# - Renames the F1/F2/F3/S1/S2 diagram titles to their new edge-notation
#   labels.
# - On the "S2" slide, swaps the start/end connection endpoints of the
#   dashed connector (it now runs Work1 -> Work2 instead of Work2 -> Work1)
#   and clears the flip flags that were compensating for the old direction.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Footer date placeholder text: 2024-07-12 -> 2024-07-20
#    (slide master + every slide layout)
# ---------------------------------------------------------------------
$oldDate = "2024-07-12"
$newDate = "2024-07-20"

for ($j = 1; $j -le $p.SlideMaster.Shapes.Count; $j++) {
    $shp = $p.SlideMaster.Shapes.Item($j)
    if ($shp.HasTextFrame -eq -1) {
        if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

for ($i = 1; $i -le $p.SlideMaster.CustomLayouts.Count; $i++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $shp = $layout.Shapes.Item($j)
        if ($shp.HasTextFrame -eq -1) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2) Slide title text updates (title placeholder is always shape 1)
# ---------------------------------------------------------------------
$p.Slides.Item(2).Shapes.Item(1).TextFrame.TextRange.Text = "W2 |> W1 > W2"
$p.Slides.Item(3).Shapes.Item(1).TextFrame.TextRange.Text = "W1 |> W2 |> W1"
$p.Slides.Item(4).Shapes.Item(1).TextFrame.TextRange.Text = "W1 <|> W2; W1 > W2"
$p.Slides.Item(5).Shapes.Item(1).TextFrame.TextRange.Text = "W1 > W2"
$p.Slides.Item(6).Shapes.Item(1).TextFrame.TextRange.Text = "W1 |> W2"

# ---------------------------------------------------------------------
# 3) Slide 6 ("S2"): flip the dashed connector's start/end connection
#    sites and drop the flipH/flipV compensation.
# ---------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$connector = $s6.Shapes.Item(2)
$rect5 = $s6.Shapes.Item(3)
$rect6 = $s6.Shapes.Item(4)

$connector.ConnectorFormat.BeginDisconnect()
$connector.ConnectorFormat.EndDisconnect()
$connector.ConnectorFormat.BeginConnect($rect5, 2)
$connector.ConnectorFormat.EndConnect($rect6, 1)

$connector.HorizontalFlip = $false
$connector.VerticalFlip = $false
